$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 2 with new data (single form values)
$ws.Range("B2").Value = "Lê Thị Yến"
$ws.Range("C2").Value = "Bắc Ninh"
$ws.Range("E2").Value = "Hà Nội"

# Update the active selection to F3 (matches saved view state in diff)
$ws.Range("F3").Select()
